$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contracts")

$ws.Cells.Item(6, 1).Value = "IAaveDistributionManager"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 1).Value = "IAaveIncentivesController"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = "IAToken"
$ws.Cells.Item(8, 2).Value = 2
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = "IERC20"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = "ILendingPool"
$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = "ILendingPoolAddressesProvider"
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = "IPriceOracleGetter"
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = "IScaledBalanceToken"
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = "IVariableDebtToken"
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = "IEntryPositionsManager"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = "IExitPositionsManager"
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = "IGetterUnderlyingAsset"
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(18, 1).Value = "IIncentivesVault"
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(19, 1).Value = "IInterestRatesManager"
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(20, 1).Value = "IMorpho"
$ws.Cells.Item(20, 2).Value = 8
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(21, 1).Value = "IOracle"
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(22, 1).Value = "IRewardsManager"
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(23, 1).Value = "ILido"
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(24, 1).Value = "IndexesLens"
$ws.Cells.Item(24, 2).Value = 3
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(25, 1).Value = "ILens"
$ws.Cells.Item(25, 2).Value = 3
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(26, 1).Value = "Lens"
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(27, 1).Value = "LensStorage"
$ws.Cells.Item(27, 2).Value = 11
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(28, 1).Value = "MarketsLens"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(29, 1).Value = "RatesLens"
$ws.Cells.Item(29, 2).Value = 2
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(30, 1).Value = "UsersLens"
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 1).Value = "DataTypes"
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(32, 1).Value = "Errors"
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(33, 1).Value = "ReserveConfiguration"
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 1).Value = "UserConfiguration"
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(35, 1).Value = "InterestRatesModel"
$ws.Cells.Item(35, 2).Value = 4
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 1).Value = "Types"
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(37, 1).Value = "MatchingEngine"
$ws.Cells.Item(37, 2).Value = 1
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(38, 1).Value = "Morpho"
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(39, 1).Value = "MorphoGovernance"
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = "MorphoStorage"
$ws.Cells.Item(40, 2).Value = 10
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(41, 1).Value = "MorphoUtils"
$ws.Cells.Item(41, 2).Value = 9
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(42, 1).Value = "PositionsManagerUtils"
$ws.Cells.Item(42, 2).Value = 3
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(47, 1).Value = "ICEth"
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(48, 1).Value = "IComptroller"
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(49, 1).Value = "IInterestRateModel"
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(50, 1).Value = "ICToken"
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(51, 1).Value = "ICEther"
$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(52, 1).Value = "ICompoundOracle"
$ws.Cells.Item(52, 2).Value = 0
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(53, 1).Value = "IIncentivesVault"
$ws.Cells.Item(53, 2).Value = 1
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(54, 1).Value = "IInterestRatesManager"
$ws.Cells.Item(54, 2).Value = 0
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(55, 1).Value = "IMorpho"
$ws.Cells.Item(55, 2).Value = 5
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(56, 1).Value = "IOracle"
$ws.Cells.Item(56, 2).Value = 0
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(57, 1).Value = "IPositionsManager"
$ws.Cells.Item(57, 2).Value = 0
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(58, 1).Value = "IRewardsManager"
$ws.Cells.Item(58, 2).Value = 1
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(59, 1).Value = "IWETH"
$ws.Cells.Item(59, 2).Value = 0
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(60, 1).Value = "IndexesLens"
$ws.Cells.Item(60, 2).Value = 2
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(61, 1).Value = "ILens"
$ws.Cells.Item(61, 2).Value = 3
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(62, 1).Value = "Lens"
$ws.Cells.Item(62, 2).Value = 1
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(63, 1).Value = "LensStorage"
$ws.Cells.Item(63, 2).Value = 9
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(64, 1).Value = "MarketsLens"
$ws.Cells.Item(64, 2).Value = 1
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(65, 1).Value = "RatesLens"
$ws.Cells.Item(65, 2).Value = 1
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(66, 1).Value = "RewardsLens"
$ws.Cells.Item(66, 2).Value = 1
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(67, 1).Value = "UsersLens"
$ws.Cells.Item(67, 2).Value = 1
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(68, 1).Value = "CompoundMath"
$ws.Cells.Item(68, 2).Value = 0
$ws.Cells.Item(68, 3).Value = 0
$ws.Cells.Item(69, 1).Value = "InterestRatesModel"
$ws.Cells.Item(69, 2).Value = 4
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(70, 1).Value = "Types"
$ws.Cells.Item(70, 2).Value = 0
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(71, 1).Value = "MatchingEngine"
$ws.Cells.Item(71, 2).Value = 1
$ws.Cells.Item(71, 3).Value = 0
$ws.Cells.Item(72, 1).Value = "Morpho"
$ws.Cells.Item(72, 2).Value = 1
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(73, 1).Value = "MorphoGovernance"
$ws.Cells.Item(73, 2).Value = 1
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(74, 1).Value = "MorphoStorage"
$ws.Cells.Item(74, 2).Value = 9
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(75, 1).Value = "MorphoUtils"
$ws.Cells.Item(75, 2).Value = 5
$ws.Cells.Item(75, 3).Value = 0
$ws.Cells.Item(76, 1).Value = "PositionsManager"
$ws.Cells.Item(76, 2).Value = 3
$ws.Cells.Item(76, 3).Value = 0
$ws.Cells.Item(77, 1).Value = "RewardsManager"
$ws.Cells.Item(77, 2).Value = 4
$ws.Cells.Item(77, 3).Value = 0
